$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old empty spacer row (row 3) is removed; rows 4-6 shift up to 3-5.
$ws.Rows(3).Delete()

# The table used to report area for three census years (1989/2002/2014) in
# columns B:D. Only the latest (2014) figure is kept, so the extra year
# columns are dropped entirely (not just cleared) -- column B stays in place.
$ws.Range("C1:D1").EntireColumn.Delete()

# Drop the now-orphaned subtitle ("according to the population census data")
# and the stray formatted placeholder cell next to the title.
$ws.Range("A2:B2").Clear()
$ws.Range("B1").Clear()

# Keep the single remaining year/value pair, updated for the simplified table.
$ws.Range("B4").Value = 2014
$ws.Range("B5").Value = 646.70000000000005

# Now that column B is the last (and only) data column, its borders close
# off the little box instead of sitting in the middle of it.
$ws.Range("B4").Borders.Item(7).LineStyle = 1
$ws.Range("B4").Borders.Item(7).Weight = 2
$ws.Range("B4").Borders.Item(10).LineStyle = 1
$ws.Range("B4").Borders.Item(10).Weight = -4138

$ws.Range("B5").Borders.Item(7).LineStyle = -4142
$ws.Range("B5").Borders.Item(10).LineStyle = 1
$ws.Range("B5").Borders.Item(10).Weight = -4138

# Add back a trailing blank spacer row (new row 6), matching the old layout.
$ws.Rows(6).Insert()
$ws.Range("A6:B6").ClearContents()
$ws.Range("A6:B6").Style = "Normal"

# All rows now use a taller, consistent custom height.
for ($r = 1; $r -le 6; $r++) {
    $ws.Rows($r).RowHeight = 20.1
}
